$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cell = $t.Cell(1, 1)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "32÷9=3, 5"

$cell = $t.Cell(1, 2)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "62÷7=8, 6"

$cell = $t.Cell(1, 3)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "50÷9=5, 5"

$cell = $t.Cell(1, 4)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "12÷7=1, 5"

$cell = $t.Cell(1, 5)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "17÷7=2, 3"

$cell = $t.Cell(5, 1)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "44÷6=7, 2"

$cell = $t.Cell(5, 2)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "19÷7=2, 5"

$cell = $t.Cell(5, 3)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "21÷9=2, 3"

$cell = $t.Cell(5, 4)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "97÷3=32, 1"

$cell = $t.Cell(5, 5)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "30÷4=7, 2"

$cell = $t.Cell(9, 1)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "77÷2=38, 1"

$cell = $t.Cell(9, 2)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "16÷6=2, 4"

$cell = $t.Cell(9, 3)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "88÷5=17, 3"

$cell = $t.Cell(9, 4)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "84÷4=21, 0"

$cell = $t.Cell(9, 5)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "26÷2=13, 0"

$cell = $t.Cell(13, 1)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "24÷2=12, 0"

$cell = $t.Cell(13, 2)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "45÷7=6, 3"

$cell = $t.Cell(13, 3)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "46÷7=6, 4"

$cell = $t.Cell(13, 4)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "64÷6=10, 4"

$cell = $t.Cell(13, 5)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "79÷6=13, 1"

$cell = $t.Cell(17, 1)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "96÷2=48, 0"

$cell = $t.Cell(17, 2)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "50÷8=6, 2"

$cell = $t.Cell(17, 3)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "71÷4=17, 3"

$cell = $t.Cell(17, 4)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "54÷7=7, 5"

$cell = $t.Cell(17, 5)
$r0 = $cell.Range
$trimmed = $r0.Text.TrimEnd([char]13, [char]7)
$newRange = $d.Range($r0.Start, $r0.Start + $trimmed.Length)
$newRange.Text = "52÷3=17, 1"
